$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 18:35"

# Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 1576007
$ws.Range("C4").Value = 5424
$ws.Range("D4").Value = 361771
$ws.Range("E4").Value = 1120389
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 93847

# Turquia (row 12) - updated case counts
$ws.Range("B12").Value = 152587
$ws.Range("C12").Value = 972
$ws.Range("D12").Value = 113987
$ws.Range("E12").Value = 34378
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 23
$ws.Range("H12").Value = 4222

# Santa Lucia / Belice swap rank positions (rows 196-197)
$ws.Range("A196").Value = "Belice"
$ws.Range("D196").Value = 16
$ws.Range("H196").Value = 2

$ws.Range("A197").Value = "Santa Lucia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# Montserrat / Groenlandia swap rank positions (rows 209-210)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
